# Fuel Prod Imp Exp Balancing Priorities.xlsx - "updated 4.0 files and mdl"
#
# Data-level changes applied (the only ones a headless Excel COM session can
# actually make - cosmetic re-save artifacts such as the theme display name,
# fileVersion/rupBuild, revisionPtr GUIDs, x14ac:dyDescent font-metric caches
# and sub-pixel "best fit" column widths are produced internally by a given
# Excel build when it resaves a file and are not exposed on the object model,
# so they are intentionally left alone here):
#
#   About!C1   : "last updated" date 2024-01-03 -> 2024-03-28 (45294 -> 45379)
#   FPIEBP!B3:D3 ("hard coal" row): production/imports/exports priorities
#                shuffled 3,2,1 -> 1,3,2
#   FPIEBP selection moved from F4 to E3 (the sheet that was active/selected
#   when the file was last saved)

$wb = $excel.ActiveWorkbook

# --- About sheet: bump the "last updated" date shown in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: re-prioritize hard coal's production/imports/exports ---
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# Leave the saved selection where the author left it (and keep FPIEBP as the
# active/selected tab, matching activeTab=1 in the original workbook).
$wsFpiebp.Range("E3").Select()
